$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.459.38'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').Value = '2.354.56'
$ws.Range('E3').Value = '  +5.99%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''312.93'
$ws.Range('E5').Value = '  +5.55%  '
$ws.Range('D6').Value = '''109.38'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.639'
$ws.Range('E9').Value = '  +6.41%  '
$ws.Range('D10').Value = '''43.38'
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').Value = '''0.0940'
$ws.Range('E11').Value = '  +2.84%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  +3.35%  '
$ws.Range('E14').Value = '  +2.40%  '
$ws.Range('D15').Value = '''16.43'
$ws.Range('E15').Value = '  +8.71%  '
$ws.Range('D16').Value = '2.707.60'
$ws.Range('E16').Value = '  +6.10%  '
$ws.Range('D17').Value = '2.419.80'
$ws.Range('E17').Value = '  +8.79%  '
$ws.Range('D18').Value = '43.435.81'
$ws.Range('E18').Value = '  +2.65%  '
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').Value = '''74.97'
$ws.Range('E21').Value = '  +3.75%  '
$ws.Range('D22').Value = '''3.45'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').Value = '''2.56'
$ws.Range('E23').Value = '  +10.50%  '
$ws.Range('D24').Value = '''258.76'
$ws.Range('E24').Value = '  +13.25%  '
$ws.Range('D25').Value = '''9.23'
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('E26').Value = '  +3.87%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''39.39'
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').Value = '''2.25'
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').Value = '''22.54'
$ws.Range('E30').Value = '  +7.38%  '
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').Value = '''173.34'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').Value = '''6.04'
$ws.Range('E34').Value = '  +7.51%  '
$ws.Range('E35').Value = '  +5.61%  '
$ws.Range('D36').Value = '''4.97'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '''4.16'
$ws.Range('E37').Value = '  -4.15%  '
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('D40').Value = '''2.75'
$ws.Range('E40').Value = '  +13.89%  '
$ws.Range('D41').Value = '''72.42'
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('E42').Value = '  +13.77%  '
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('D44').Value = '''12.81'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  +4.08%  '
$ws.Range('D47').Value = '''9.38'
$ws.Range('E47').Value = '  +11.32%  '
$ws.Range('D48').Value = '''110.97'
$ws.Range('E48').Value = '  +7.55%  '
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').Value = '''0.100'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('D51').Value = '''0.467'
$ws.Range('E51').Value = '  +6.73%  '
